$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (Volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/2/2024  Through  9/8/2024"

# --- Cells that change from a numeric value to the text placeholder "0" (shared string 20) ---
# Template source cells (style 14, text "0"): C14, D14, F14, G14 (untouched by this edit)
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))

# --- Cells that change from a numeric value to the text placeholder "***.*" (shared string 21) ---
# Template source cells (style 14, text "***.*"): E14, H14 (untouched by this edit)
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("E14").Copy($ws.Range("E27"))

# --- Cells that change from text placeholder to a genuine numeric value ---
# Use a stable numeric-style template cell first (to fix style/type), then set the real value
$ws.Range("I15").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 2
$ws.Range("K15").Copy($ws.Range("E28"))
$ws.Range("E28").Value = 0

# --- Plain numeric value updates ---
$ws.Range("M14").Value = -83.333333333333
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = -47.826086956521
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -7.142857142857
$ws.Range("I16").Value = 122
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 1.666666666666
$ws.Range("L16").Value = -4.6875
$ws.Range("M16").Value = -35.789473684210
$ws.Range("N16").Value = -83.905013192612
$ws.Range("C17").Value = 6
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 20
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 255
$ws.Range("J17").Value = 258
$ws.Range("K17").Value = -1.162790697674
$ws.Range("L17").Value = -6.25
$ws.Range("M17").Value = 65.584415584415
$ws.Range("N17").Value = -53.720508166969
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 6
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 60
$ws.Range("J18").Value = 76
$ws.Range("K18").Value = -21.052631578947
$ws.Range("L18").Value = -45.454545454545
$ws.Range("M18").Value = -49.152542372881
$ws.Range("N18").Value = -91.279069767441
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 3.030303030303
$ws.Range("I19").Value = 249
$ws.Range("J19").Value = 340
$ws.Range("K19").Value = -26.764705882352
$ws.Range("L19").Value = -31.967213114754
$ws.Range("M19").Value = -41.822429906542
$ws.Range("N19").Value = -48.659793814433
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -36.363636363636
$ws.Range("J20").Value = 75
$ws.Range("K20").Value = 1.333333333333
$ws.Range("L20").Value = -11.627906976744
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -90.488110137672
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -4.545454545454
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = -14.736842105263
$ws.Range("I21").Value = 775
$ws.Range("J21").Value = 880
$ws.Range("K21").Value = -11.931818181818
$ws.Range("L21").Value = -20.349434737923
$ws.Range("M21").Value = -21.079429735234
$ws.Range("N21").Value = -76.670680313064
$ws.Range("F22").Value = 3
$ws.Range("I22").Value = 21
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = 5
$ws.Range("M22").Value = -16
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = -11.111111111111
$ws.Range("I23").Value = 95
$ws.Range("J23").Value = 115
$ws.Range("K23").Value = -17.391304347826
$ws.Range("L23").Value = -16.666666666666
$ws.Range("M23").Value = 11.764705882352
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -61.904761904761
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 6.849315068493
$ws.Range("I24").Value = 689
$ws.Range("J24").Value = 713
$ws.Range("K24").Value = -3.366058906030
$ws.Range("L24").Value = -19.131455399061
$ws.Range("M24").Value = -22.497187851518
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 7.142857142857
$ws.Range("I25").Value = 191
$ws.Range("J25").Value = 229
$ws.Range("K25").Value = -16.593886462882
$ws.Range("L25").Value = -38.782051282051
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = -3.571428571428
$ws.Range("I26").Value = 528
$ws.Range("J26").Value = 436
$ws.Range("K26").Value = 21.100917431192
$ws.Range("L26").Value = 32.663316582914
$ws.Range("M26").Value = 41.176470588235
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 39
$ws.Range("K28").Value = 2.564102564102
$ws.Range("L28").Value = 8.108108108108
$ws.Range("M29").Value = -93.333333333333
$ws.Range("N29").Value = -97.916666666666
$ws.Range("M30").Value = -91.666666666666
$ws.Range("N30").Value = -97.368421052631
